# Update "想去人数" (F column) figures across all sheets to reflect
# the latest scrape snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 2697
$ws1.Range("F3").Value = 1042
$ws1.Range("F4").Value = 19293
$ws1.Range("F6").Value = 2182
$ws1.Range("F7").Value = 736
$ws1.Range("F9").Value = 416
$ws1.Range("F10").Value = 672
$ws1.Range("F12").Value = 240
$ws1.Range("F14").Value = 351
$ws1.Range("F16").Value = 252
$ws1.Range("F18").Value = 173
$ws1.Range("F19").Value = 16
$ws1.Range("F20").Value = 18

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 187
$ws2.Range("F6").Value = 119
$ws2.Range("F7").Value = 270
$ws2.Range("F8").Value = 124

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5971
$ws3.Range("F3").Value = 628
$ws3.Range("F4").Value = 580

# Sheet "全部类型" (All Types) - aggregated view of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5971
$ws4.Range("F3").Value = 628
$ws4.Range("F4").Value = 580
$ws4.Range("F5").Value = 187
$ws4.Range("F7").Value = 2697
$ws4.Range("F8").Value = 1042
$ws4.Range("F9").Value = 19293
$ws4.Range("F13").Value = 119
$ws4.Range("F14").Value = 270
$ws4.Range("F15").Value = 2182
$ws4.Range("F16").Value = 736
$ws4.Range("F17").Value = 124
$ws4.Range("F19").Value = 416
$ws4.Range("F20").Value = 672
$ws4.Range("F22").Value = 240
$ws4.Range("F27").Value = 351
$ws4.Range("F30").Value = 252
$ws4.Range("F34").Value = 173
$ws4.Range("F37").Value = 16
$ws4.Range("F40").Value = 18
